# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# 1. Insert a new "Player Info" sheet as the first sheet, with the
#    player's basic info (ID/NAME/BATTING_HAND/BOWL_STYLE).
# 2. On "ODI Batting" and "ODI Bowling", rename the MATCH_CARD_LINK
#    column to MATCH_CODE, and replace the full scorecard URL values
#    with just the numeric match code that used to be the query string.

$wb = $excel.ActiveWorkbook

$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

# --- ODI Batting: MATCH_CARD_LINK (col D) -> MATCH_CODE -------------------
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "'4401"
$batting.Range("D3").Value = "'4405"
$batting.Range("D4").Value = "'4408"
$batting.Range("D5").Value = "'4421"
$batting.Range("D6").Value = "'4460"

# --- ODI Bowling: MATCH_CARD_LINK (col B) -> MATCH_CODE -------------------
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'4401"
$bowling.Range("B3").Value = "'4408"
$bowling.Range("B4").Value = "'4421"
$bowling.Range("B5").Value = "'4460"

# --- New "Player Info" sheet, inserted before "ODI Batting" ---------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$playerInfo.Range("A2").Value = "'4849"
$playerInfo.Range("B2").Value = "Lubabalo Lutho Sipamla"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"
